$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 697
$ws.Range("J3").Value = 769
$ws.Range("J4").Value = 167
$ws.Range("J5").Value = 55
$ws.Range("I6").Value = 8966
$ws.Range("J6").Value = 1121
$ws.Range("I7").Value = 26183
$ws.Range("J7").Value = 2809

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J6").Value = 24
$ws.Range("J7").Value = 83
$ws.Range("J8").Value = 181
$ws.Range("J9").Value = 19
$ws.Range("J10").Value = 16
$ws.Range("J11").Value = 38
$ws.Range("J15").Value = 33
$ws.Range("J19").Value = 90
$ws.Range("J20").Value = 55
$ws.Range("J21").Value = 4
$ws.Range("J23").Value = 25
$ws.Range("J25").Value = 17
$ws.Range("J29").Value = 140
$ws.Range("J31").Value = 21
$ws.Range("J33").Value = 117
$ws.Range("J34").Value = 17
$ws.Range("J35").Value = 5
$ws.Range("J36").Value = 42
$ws.Range("J37").Value = 103
$ws.Range("J41").Value = 19
$ws.Range("J42").Value = 127
$ws.Range("J44").Value = 24
$ws.Range("I48").Value = 329
$ws.Range("J51").Value = 37
$ws.Range("J52").Value = 64
$ws.Range("J53").Value = 29
$ws.Range("J54").Value = 48
$ws.Range("J55").Value = 33
$ws.Range("J64").Value = 18
$ws.Range("J65").Value = 67
$ws.Range("J67").Value = 105
$ws.Range("J70").Value = 6
$ws.Range("J75").Value = 12
$ws.Range("J76").Value = 45
$ws.Range("J77").Value = 24
$ws.Range("J84").Value = 34
$ws.Range("J85").Value = 114
$ws.Range("J89").Value = 35
$ws.Range("J91").Value = 35
$ws.Range("J92").Value = 7
$ws.Range("J98").Value = 20
$ws.Range("I101").Value = 26183
$ws.Range("J101").Value = 2809

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J5").Value = 4
$ws.Range("J7").Value = 114

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J6").Value = 26
$ws.Range("J7").Value = 64

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 11
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 57
$ws.Range("J3").Value = 58
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 181

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 29

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 83

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J6").Value = 36
$ws.Range("J7").Value = 103

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 105

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J6").Value = 2
$ws.Range("J7").Value = 21

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 34

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J3").Value = 18
$ws.Range("J6").Value = 25
$ws.Range("J7").Value = 67

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 30
$ws.Range("J6").Value = 50
$ws.Range("J7").Value = 117

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J2").Value = 14
$ws.Range("J6").Value = 24
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 45
$ws.Range("J4").Value = 4
$ws.Range("J7").Value = 140

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 28
$ws.Range("J6").Value = 40
$ws.Range("J7").Value = 90

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 24

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I6").Value = 169
$ws.Range("I7").Value = 329

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J3").Value = 7
$ws.Range("J6").Value = 11
$ws.Range("J7").Value = 24

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J2").Value = 4
$ws.Range("J7").Value = 19

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 18
$ws.Range("J6").Value = 85
$ws.Range("J7").Value = 127

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J2").Value = 6
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J2").Value = 7
$ws.Range("J7").Value = 25

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("J6").Value = 3
$ws.Range("J7").Value = 4

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 18

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 19
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 55

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 12
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 17

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 17

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J3").Value = 10
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 20

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("J6").Value = 4
$ws.Range("J7").Value = 5

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J3").Value = 6
$ws.Range("J7").Value = 19

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("J6").Value = 1
$ws.Range("J7").Value = 7

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("J2").Value = 3
$ws.Range("J7").Value = 6

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J6").Value = 3
$ws.Range("J7").Value = 12

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 11
$ws.Range("J3").Value = 14
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 24
